$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5938.294
$ws.Range("I40").Value = 4717.222
$ws.Range("K40").Value = 4717.222
$ws.Range("M40").Value = -4542.222

$ws.Range("H58").Value = 314
$ws.Range("I58").Value = 314
$ws.Range("K58").Value = 942
$ws.Range("M58").Value = -792

$ws.Range("H98").Value = 729.86957
$ws.Range("I98").Value = 729.86957
$ws.Range("K98").Value = 729.86957
$ws.Range("M98").Value = 768.13043

$ws.Range("H112").Value = 10871821
$ws.Range("J112").Value = 11906932
$ws.Range("L112").Value = 35720796
$ws.Range("N112").Value = -35723012

$ws.Range("H122").Value = 729.86957
$ws.Range("I122").Value = 729.86957
$ws.Range("K122").Value = 2189.60871
$ws.Range("M122").Value = 260.39129

$ws.Range("H137").Value = 7817.067
$ws.Range("I137").Value = 3717.375
$ws.Range("J137").Value = 12502.429
$ws.Range("K137").Value = 11152.125
$ws.Range("L137").Value = 37507.287
$ws.Range("M137").Value = -8602.125
$ws.Range("N137").Value = -42607.287

$ws.Range("H138").Value = 58825820
$ws.Range("I138").Value = 1656.2727
$ws.Range("K138").Value = 4968.8181
$ws.Range("M138").Value = 171.1818999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13339840
$ws.Range("I32").Value = 15628920
$ws.Range("K32").Value = 15628920
$ws.Range("M32").Value = -15628633

$ws.Range("H41").Value = 21203.8
$ws.Range("J41").Value = 39000
$ws.Range("L41").Value = 39000
$ws.Range("N41").Value = -39828

$ws.Range("H61").Value = 25002974
$ws.Range("I61").Value = 31251842
$ws.Range("K61").Value = 31251842
$ws.Range("M61").Value = -31251630

$ws.Range("H74").Value = 50057464
$ws.Range("I74").Value = 77010400
$ws.Range("J74").Value = 2002.1428
$ws.Range("K74").Value = 77010400
$ws.Range("L74").Value = 2002.1428
$ws.Range("M74").Value = -77009526
$ws.Range("N74").Value = -3750.1428

$ws.Range("H77").Value = 50057464
$ws.Range("I77").Value = 77010400
$ws.Range("J77").Value = 2002.1428
$ws.Range("K77").Value = 385052000
$ws.Range("L77").Value = 10010.714
$ws.Range("M77").Value = -385047632
$ws.Range("N77").Value = -18746.714

$ws.Range("H125").Value = 60999.332
$ws.Range("J125").Value = 66499
$ws.Range("L125").Value = 66499
$ws.Range("N125").Value = -76339

$ws.Range("H132").Value = 26389506
$ws.Range("I132").Value = 11978.129
$ws.Range("K132").Value = 35934.387
$ws.Range("M132").Value = -33404.387

$ws.Range("H136").Value = 25002974
$ws.Range("I136").Value = 31251842
$ws.Range("K136").Value = 93755526
$ws.Range("M136").Value = -93752976

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 90
$ws.Range("I22").Value = 90
$ws.Range("K22").Value = 90
$ws.Range("M22").Value = 83

$ws.Range("H80").Value = 986.44446
$ws.Range("J80").Value = 646.8889
$ws.Range("L80").Value = 646.8889
$ws.Range("N80").Value = -2642.8889

$ws.Range("H83").Value = 986.44446
$ws.Range("J83").Value = 646.8889
$ws.Range("L83").Value = 3234.4445
$ws.Range("N83").Value = -13218.4445

$ws.Range("H117").Value = 39000
$ws.Range("J117").Value = 39000
$ws.Range("L117").Value = 39000
$ws.Range("N117").Value = -48178

$ws.Range("H124").Value = 75516.664
$ws.Range("J124").Value = 75516.664
$ws.Range("L124").Value = 75516.664
$ws.Range("N124").Value = -85336.664

$ws.Range("H135").Value = 99999
$ws.Range("J135").Value = 99999
$ws.Range("L135").Value = 99999
$ws.Range("N135").Value = -110139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 76927860
$ws.Range("I31").Value = 1413.5
$ws.Range("K31").Value = 1413.5
$ws.Range("M31").Value = -1118.5

$ws.Range("H34").Value = 76927860
$ws.Range("I34").Value = 1413.5
$ws.Range("K34").Value = 1413.5
$ws.Range("M34").Value = -1211.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2324
$ws.Range("J98").Value = 2324
$ws.Range("L98").Value = 6972
$ws.Range("N98").Value = -9968

$ws.Range("H140").Value = 1754.85
$ws.Range("I140").Value = 770.6923
$ws.Range("K140").Value = 2312.0769
$ws.Range("M140").Value = 2867.9231

$ws.Range("H141").Value = 4290.524
$ws.Range("I141").Value = 1269.2858
$ws.Range("K141").Value = 3807.8574
$ws.Range("M141").Value = 1372.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 578
$ws.Range("I107").Value = 506.8
$ws.Range("K107").Value = 506.8
$ws.Range("M107").Value = 1413.2

$ws.Range("H126").Value = 9408931
$ws.Range("I126").Value = 4593535
$ws.Range("K126").Value = 13780605
$ws.Range("M126").Value = -13778135

$ws.Range("H132").Value = 3592.6667
$ws.Range("I132").Value = 3662.3
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 10986.9
$ws.Range("L132").Value = 6600
$ws.Range("M132").Value = -8456.900000000001
$ws.Range("N132").Value = -11660

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 968.8570999999999
$ws.Range("J55").Value = 1463
$ws.Range("L55").Value = 1463
$ws.Range("N55").Value = -1809

$ws.Range("H100").Value = 4390
$ws.Range("I100").Value = 3425.5557
$ws.Range("J100").Value = 5475
$ws.Range("K100").Value = 3425.5557
$ws.Range("L100").Value = 5475
$ws.Range("M100").Value = -2884.5557
$ws.Range("N100").Value = -6557

$ws.Range("H127").Value = 94357
$ws.Range("J127").Value = 94357
$ws.Range("L127").Value = 94357
$ws.Range("N127").Value = -104277

$ws.Range("H132").Value = 64523172
$ws.Range("I132").Value = 4876.3335
$ws.Range("K132").Value = 14629.0005
$ws.Range("M132").Value = -12099.0005

$ws.Range("H136").Value = 2381.7083
$ws.Range("I136").Value = 1984.9333
$ws.Range("K136").Value = 5954.7999
$ws.Range("M136").Value = -3404.7999

$ws.Range("H137").Value = 69497.5
$ws.Range("I137").Value = 40000
$ws.Range("K137").Value = 40000
$ws.Range("M137").Value = -34900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 142864620
$ws.Range("J41").Value = 10039.2
$ws.Range("L41").Value = 10039.2
$ws.Range("N41").Value = -10819.2

$ws.Range("H100").Value = 28858692
$ws.Range("I100").Value = 34829096
$ws.Range("K100").Value = 69658192
$ws.Range("M100").Value = -69657651

$ws.Range("H107").Value = 2392.8235
$ws.Range("I107").Value = 1888.4445
$ws.Range("J107").Value = 2960.25
$ws.Range("K107").Value = 5665.333500000001
$ws.Range("L107").Value = 8880.75
$ws.Range("M107").Value = -3745.333500000001
$ws.Range("N107").Value = -12720.75

$ws.Range("H122").Value = 71502510
$ws.Range("I122").Value = 91001016
$ws.Range("J122").Value = 7997
$ws.Range("K122").Value = 273003048
$ws.Range("L122").Value = 23991
$ws.Range("M122").Value = -273000598
$ws.Range("N122").Value = -28891

$ws.Range("H128").Value = 149999
$ws.Range("J128").Value = 149999
$ws.Range("L128").Value = 149999
$ws.Range("N128").Value = -159959

$ws.Range("H132").Value = 2364.4915
$ws.Range("I132").Value = 2242.4583
$ws.Range("J132").Value = 2897
$ws.Range("K132").Value = 6727.374899999999
$ws.Range("L132").Value = 8691
$ws.Range("M132").Value = -4197.374899999999
$ws.Range("N132").Value = -13751

$ws.Range("H136").Value = 2150.923
$ws.Range("I136").Value = 1055.3529
$ws.Range("K136").Value = 3166.0587
$ws.Range("M136").Value = -616.0587000000005
